$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")
$ws.Range("A1").Value = "Test"
